$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New BOM data (designator, MPN, qty) for rows 2..20.
# Row 9 (IC2,IC4) keeps its Link hyperlink text; all others clear Link.
$data = @(
    @("C11,C17", "GCM21BR71H474KA55L", 4, ""),
    @("C4,C5,C6,C8,C9,C10,C12,C14,C15,C16,C18", "GRM21BR71H104KA01L", 22, ""),
    @("C2,C7,C13,C19", "GRM21BR61H475KE51L", 8, ""),
    @("D3,D4,D8,D9,D10,D11", "150060RS75000", 12, ""),
    @("D1,D2", "150060VS75000", 4, ""),
    @("D6,D7", "CDBU0340", 4, ""),
    @("IC1,IC3", "SI8650AB-B-IS1", 4, ""),
    @("IC2,IC4", "A3921KLPTR-T", 4, "https://www.digikey.jp/product-detail/ja/allegro-microsystems/A3921KLPTR-T/620-1523-1-ND/4318335"),
    @("Q1,Q2,Q3,Q4", "NVMFD5C466NLT1G", 8, ""),
    @("R10,R30,R31,R32,R33", "CRCW060310K0JNEAC", 10, ""),
    @("R24,R27", "RC0603FR-073KL", 4, ""),
    @("R25,R28", "CRCW06033K74FKEA", 4, ""),
    @("R26,R29", "ERJ-3RBD1201V", 4, ""),
    @("R4,R5,R12,R13,R16,R17,R20,R21", "ERJ-3RED21R0V", 16, ""),
    @("U1,U4", "AZ2085D-ADJTRG1", 4, ""),
    @("U2", "STM32F042K6T6TR", 2, ""),
    @("U3", "MCP2562-E/SN", 2, ""),
    @("C1,C3", "CL21B105KBFNNNE", 4, ""),
    @("RN1", "EXB-38V222JV", 2, "")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    if ($entry[3] -ne "") {
        $ws.Cells.Item($row, 4).Value = $entry[3]
    } else {
        $ws.Cells.Item($row, 4).Value = $null
    }
    $row = $row + 1
}

# Clear the two now-unused trailing blank rows (previously rows 23 & 24).
$ws.Range("A23:E24").Clear()

# Left-align designator/MPN/qty columns for the data rows, matching new styling.
$ws.Range("A2:D18").HorizontalAlignment = -4131

# The last two BOM rows (new parts) use a slightly different left-aligned
# style (no forced vertical centering) so they end up on a distinct xf.
$ws.Range("A19:C20").HorizontalAlignment = -4131
$ws.Range("A19:C20").VerticalAlignment = -4107
$ws.Range("D19:D20").HorizontalAlignment = -4131

# Update selection to match new active cell.
$ws.Range("A45").Select()
